# Weekly refresh of Fruta/Hortaliza data: rows 2-40 get their
# Fecha/Volumen/Precio columns (D, J, K, L, M, P) re-shuffled to reflect
# the latest weekly pull. Row 30 is unaffected (its price-report record
# didn't move this week).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns whose values travel together as one "record" when rows are
# re-ordered by the weekly import.
$cols = @("D", "J", "K", "L", "M", "P")

# Destination row -> source row (values currently sitting in $map[$dest]
# before this edit are the ones that should end up in row $dest).
$map = @{
    2 = 12
    3 = 16
    4 = 24
    5 = 37
    6 = 39
    7 = 4
    8 = 6
    9 = 13
    10 = 21
    11 = 7
    12 = 14
    13 = 20
    14 = 34
    15 = 8
    16 = 18
    17 = 15
    18 = 40
    19 = 35
    20 = 19
    21 = 31
    22 = 26
    23 = 5
    24 = 32
    25 = 36
    26 = 23
    27 = 11
    28 = 3
    29 = 10
    30 = 30
    31 = 38
    32 = 27
    33 = 28
    34 = 22
    35 = 2
    36 = 33
    37 = 25
    38 = 29
    39 = 17
    40 = 9
}

# Snapshot every source cell first so that writes to earlier destination
# rows don't clobber values still needed for later destination rows.
$snapshot = @{}
foreach ($r in 2..40) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Range("$c$r").Value2
    }
    $snapshot[$r] = $rowVals
}

foreach ($dest in $map.Keys) {
    $src = $map[$dest]
    $rowVals = $snapshot[$src]
    foreach ($c in $cols) {
        $ws.Range("$c$dest").Value = $rowVals[$c]
    }
}
